# Refactoring for final DIGS data set
# Adds a new data row (ehbv-meta.6-sphenodon / Sphenodon) to the bottom of
# the table on Sheet1, mirroring the formatting of the preceding row, and
# updates the sheet selection to reflect where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append row 61, cloning the formatting (styles) of row 60 -------------
$ws.Range("A60:P60").Copy()
$ws.Range("A61:P61").PasteSpecial(-4122)   # xlPasteFormats

# --- Populate the new row's values -----------------------------------------
$ws.Range("A61").Value = "ehbv-meta.6-sphenodon"
$ws.Range("B61").Value = "meta.6-sphenodon"
$ws.Range("C61").Value = "Metahepadnavirus"
$ws.Range("D61").Value = "Endogenous metahepadnavirus 1"
$ws.Range("E61").Value = "ehbv-meta.6-sphenodon"
$ws.Range("F61").Value = 1
$ws.Range("G61").Value = "fragment"
$ws.Range("H61").Value = "nd"
$ws.Range("I61").Value = "nd"
$ws.Range("J61").Value = "nd"
$ws.Range("K61").Value = "nd"
$ws.Range("L61").Value = "Hepadnaviridae"
# M61 intentionally left blank (matches M60)
$ws.Range("N61").Value = "Sphenodon"
$ws.Range("O61").Value = "N/A"
$ws.Range("P61").Value = "ND"

# --- Restore the author's last on-screen selection --------------------------
$ws.Range("D56").Select()

Write-Output "Row 61 added; workbook updated."
